$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "batch"
$ws.Range("C2").Value = "13-02-2020"

$ws.Columns.Item(3).AutoFit() | Out-Null

$ws.Range("C1").Select() | Out-Null
